$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.866.93'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '1.563.01'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'205.66"
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'21.78"
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('E9').Value = '  -0.20%  '
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').Value = '1.786.05'
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('D13').Value = '1.566.23'
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').Value = "'0.516"
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = '26.878.66'
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').Value = "'61.27"
$ws.Range('E17').Value = '  -2.56%  '
$ws.Range('D18').Value = "'215.44"
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('D19').Value = "'7.39"
$ws.Range('E19').Value = '  +2.26%  '
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('E23').Value = '  -1.03%  '
$ws.Range('E24').Value = '  +0.97%  '
$ws.Range('D25').Value = "'154.24"
$ws.Range('E25').Value = '  +1.69%  '
$ws.Range('E26').Value = '  +1.39%  '
$ws.Range('D27').Value = "'14.96"
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('E29').Value = '  -0.62%  '
$ws.Range('E30').Value = '  +0.93%  '
$ws.Range('E31').Value = '  -3.56%  '
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('D33').Value = '1.395.46'
$ws.Range('E33').Value = '  +0.88%  '
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('E36').Value = '  -0.56%  '
$ws.Range('D37').Value = "'0.922"
$ws.Range('E37').Value = '  -2.25%  '
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  +3.25%  '
$ws.Range('D40').Value = "'0.813"
$ws.Range('E40').Value = '  +0.39%  '
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').Value = "'0.991"
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('E43').Value = '  +5.12%  '
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').Value = "'2.18"
$ws.Range('E45').Value = '  +1.00%  '
$ws.Range('D46').Value = "'63.72"
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('D47').Value = '1.699.36'
$ws.Range('D48').Value = "'86.72"
$ws.Range('E48').Value = '  +1.54%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = "'0.0504"
$ws.Range('E49').Value = '  +2.69%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0983'
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('D51').Value = "'0.0952"
$ws.Range('E51').Value = '  +1.06%  '
